# Trade #12 closed at 2026-02-17 08:13:55 - unknown UNKNOWN +0.000%
# Updates the open Trade #12 (row 13 on "All Trades" and "MarketMaking" sheets)
# to CLOSED, and rolls the resulting P&L through the Summary and Strategy
# Status sheets.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.64   # Current Capital
$summary.Range("B4").Value = -0.36     # Total P&L $
$summary.Range("B5").Value = -0.6      # Total P&L %
$summary.Range("B6").Value = 12        # Total Trades
$summary.Range("B8").Value = 6         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.64      # Capital
$status.Range("D4").Value = 12         # Trades
$status.Range("E4").Value = -0.36      # P&L $
$status.Range("F4").Value = -0.36      # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---- Row 13 (Trade #12) on both "All Trades" and "MarketMaking" sheets ----
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G13").Value = 0.35          # Exit Price
    $ws.Range("H13").Value = "CLOSED"      # Status
    $ws.Range("I13").Value = -33.9623      # P&L %
    $ws.Range("J13").Value = -0.18         # P&L $
    $ws.Range("K13").Value = 99.64         # Capital After
    $ws.Range("P13").Value = "early_exit"  # Exit Reason
    $ws.Range("Q13").Value = 5.09          # Duration (min)
}
